$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply the text-format number format (same as the header/Invoice Number column)
# to the invoice number cells A2:A7
$ws.Range("A2:A7").NumberFormat = $ws.Range("A1").NumberFormat

# Clear the Vendor column values (B2:B7) - "Burlington Hydro" entries removed
$ws.Range("B2:B7").ClearContents()

# Update the selected cell/range to match the saved selection state
$ws.Range("D7").Select()

Write-Host "done"
